$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "2+3=5"
